$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: "[Org (afko)]" table header renamed to "[Organizations]" ---
$ws.Range("A1").Value = "[Organizations]"
$ws.Range("B1").Value = "orgRef"

# --- Row 2: column headers for Organization ---
$ws.Range("B2").Value = "OrgRef"

# --- Row 8: "[PersonRegistration]" table renamed to "[Persons]", gains MiddleName column ---
$ws.Range("A8").Value = "[Persons]"
$ws.Range("B8").Value = "personRef"
$ws.Range("C8").Value = "personFirstName"
$ws.Range("D8").Value = "personMiddleName"
$ws.Range("E8").Value = "personLastName"

# --- Row 9: column headers for Person ---
$ws.Range("B9").Value = "PersonRef"
$ws.Range("C9").Value = "FirstName"
$ws.Range("D9").Value = "MiddleName"
$ws.Range("E9").Value = "LastName"

# --- Rows 10-13: Person rows lose the Ref (B) column, FirstName shifts to C, new LastName in E ---
$ws.Range("B10").Clear() | Out-Null
$ws.Range("C10").Value = "Peter"
$ws.Range("E10").Value = "Osterijen"

$ws.Range("B11").Clear() | Out-Null
$ws.Range("C11").Value = "Daniel"
$ws.Range("E11").Value = "Hoog Lieverdink"

$ws.Range("B12").Clear() | Out-Null
$ws.Range("C12").Value = "Doris"
$ws.Range("E12").Value = "Pieters-Davids"

$ws.Range("B13").Clear() | Out-Null
$ws.Range("C13").Value = "Tinus"
$ws.Range("E13").Value = "Nieuw Tonnenberg"

# --- Row 15: accPerson/accOrg headers renamed to accActor/accParty ---
$ws.Range("D15").Value = "accActor"
$ws.Range("E15").Value = "accParty"

# --- Row 16: column headers UserID->Userid, Person->Actor, Organization->Party ---
$ws.Range("B16").Value = "Userid"
$ws.Range("D16").Value = "Actor"
$ws.Range("E16").Value = "Party"

# --- Update the current selection from F16 to C9 ---
$ws.Range("C9").Select() | Out-Null
